$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F38").Value = 422
$ws.Range("G38").Value = 15365.02
$ws.Range("F47").Value = 47
$ws.Range("G47").Value = 1188.63
$ws.Range("F51").Value = 135
$ws.Range("G51").Value = 12627.9
$ws.Range("F52").Value = 22
$ws.Range("G52").Value = 1298
$ws.Range("F55").Value = 115
$ws.Range("G55").Value = 6412.4
$ws.Range("F61").Value = 208
$ws.Range("G61").Value = 54231.84
$ws.Range("F64").Value = 59
$ws.Range("G64").Value = 4691.68
$ws.Range("B66").Value = 190479.25
$ws.Range("F116").Value = 40
$ws.Range("G116").Value = 2486.8
$ws.Range("B123").Value = 70359.3
$ws.Range("B126").Value = 64196
$ws.Range("F126").Value = 1
$ws.Range("G126").Value = 32143.58
$ws.Range("B127").Value = 65258
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("F151").Value = 24
$ws.Range("G151").Value = 3200.64
$ws.Range("B155").Value = 35147.59
$ws.Range("F180").Value = 40
$ws.Range("G180").Value = 6737.6
$ws.Range("F192").Value = 21
$ws.Range("G192").Value = 1257.69
$ws.Range("B193").Value = 61796.93
$ws.Range("F206").Value = 57
$ws.Range("G206").Value = 3693.6
$ws.Range("B208").Value = 3693.6
$ws.Range("F210").Value = 103
$ws.Range("G210").Value = 5604.23
$ws.Range("B218").Value = 72094.85000000001
$ws.Range("F222").Value = 622
$ws.Range("G222").Value = 11507
$ws.Range("B229").Value = 20688.4
$ws.Range("B290").Value = 64983
$ws.Range("C290").Value = "HIM-TOTAL CARE BABY PANTS DIAPERS-M-9S"
$ws.Range("F290").Value = 6
$ws.Range("G290").Value = 514.08
$ws.Range("B291").Value = 66194
$ws.Range("C291").Value = "HIM-Total Care Baby Pants Diapers-M-9s"
$ws.Range("F291").Value = 22
$ws.Range("G291").Value = 1884.96
$ws.Range("F302").Value = 35
$ws.Range("G302").Value = 3998.05
$ws.Range("B304").Value = 63520
$ws.Range("E304").Value = 153.4
$ws.Range("F304").Value = 36
$ws.Range("G304").Value = 5194.08
$ws.Range("B305").Value = 55373
$ws.Range("E305").Value = 163.62
$ws.Range("F305").Value = -94
$ws.Range("G305").Value = -13562.32
$ws.Range("B328").Value = -17202.83
$ws.Range("F358").Value = 39
$ws.Range("G358").Value = 8978.58
$ws.Range("F359").Value = 48
$ws.Range("G359").Value = 11524.8
$ws.Range("B363").Value = 69825.69
$ws.Range("B381").Value = 47097
$ws.Range("D381").Value = 112.28
$ws.Range("E381").Value = 134.16
$ws.Range("F381").Value = 15
$ws.Range("G381").Value = 1684.2
$ws.Range("B382").Value = 58047
$ws.Range("D382").Value = 105.54
$ws.Range("E382").Value = 126.1
$ws.Range("F382").Value = 32
$ws.Range("G382").Value = 3377.28
$ws.Range("F387").Value = 400
$ws.Range("G387").Value = 38640
$ws.Range("B389").Value = 55412.14
$ws.Range("F397").Value = 69
$ws.Range("G397").Value = 2495.73
$ws.Range("F403").Value = 67
$ws.Range("G403").Value = 2716.18
$ws.Range("F404").Value = 3
$ws.Range("G404").Value = 1040.88
$ws.Range("F408").Value = 192
$ws.Range("G408").Value = 3043.2
$ws.Range("F412").Value = 80
$ws.Range("G412").Value = 997.6
$ws.Range("F416").Value = 58
$ws.Range("G416").Value = 1702.88
$ws.Range("B417").Value = 163479.88
$ws.Range("F421").Value = 53
$ws.Range("G421").Value = 4066.69
$ws.Range("B427").Value = 23185.94
$ws.Range("F453").Value = 30
$ws.Range("G453").Value = 4359.6
$ws.Range("B458").Value = 89244.06
$ws.Range("B506").Value = 60022
$ws.Range("E506").Value = 37.22
$ws.Range("F506").Value = -113
$ws.Range("G506").Value = -3709.79
$ws.Range("B507").Value = 64830
$ws.Range("E507").Value = 34.9
$ws.Range("F507").Value = 84
$ws.Range("G507").Value = 2757.72
$ws.Range("F511").Value = 204
$ws.Range("G511").Value = 20373.48
$ws.Range("F520").Value = 7
$ws.Range("G520").Value = 191.8
$ws.Range("F523").Value = 140
$ws.Range("G523").Value = 11985.4
$ws.Range("B525").Value = 115313.33
$ws.Range("F531").Value = 207
$ws.Range("G531").Value = 6853.77
$ws.Range("B535").Value = 21960.85
$ws.Range("F544").Value = 35
$ws.Range("G544").Value = 2166.5
$ws.Range("B556").Value = 41102.72
$ws.Range("F559").Value = 17
$ws.Range("G559").Value = 1766.98
$ws.Range("F560").Value = 21
$ws.Range("G560").Value = 1687.98
$ws.Range("B561").Value = 24047.61
$ws.Range("F570").Value = 4
$ws.Range("G570").Value = 2139.52
$ws.Range("B573").Value = 15428.85
$ws.Range("F610").Value = 60
$ws.Range("G610").Value = 1519.8
$ws.Range("F612").Value = 226
$ws.Range("G612").Value = 33992.66
$ws.Range("F617").Value = 8
$ws.Range("G617").Value = 384.96
$ws.Range("F620").Value = 350
$ws.Range("G620").Value = 27506.5
$ws.Range("B628").Value = 201186.31
$ws.Range("F660").Value = 46
$ws.Range("G660").Value = 1368.04
$ws.Range("F663").Value = 40
$ws.Range("G663").Value = 1430.8
$ws.Range("B668").Value = 11039.55
$ws.Range("F674").Value = 626
$ws.Range("G674").Value = 102106.86
$ws.Range("B680").Value = 103119.41
$ws.Range("F709").Value = 4
$ws.Range("G709").Value = 1237.68
$ws.Range("B713").Value = 61553.87
$ws.Range("B718").Value = 2470338.47
$ws.Range("B719").Value = 2470338.47
